$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Mã vật tư" column (column D) entirely - its header and data are
# no longer part of the import template. This shifts the old "Giá nhập"
# column (E) left into D.
$ws.Range("D1:D2").EntireColumn.Delete()

# Update the selected cell to match the saved state after the edit.
$ws.Range("H4").Select()
